# Update cryptocurrency price/volume data per the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks like a plain number (e.g. "1.00", "6.98")
# while keeping the cell's underlying type as text, matching the source
# data's inlineStr representation (no leftover number formatting is kept
# on the cell afterwards).
function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "64.703.23"
$ws.Range("E2").Value = "  -3.37%  "
$ws.Range("D3").Value = "3.406.53"
$ws.Range("E3").Value = "  -4.20%  "
$ws.Range("E4").Value = "  +0.11%  "
Set-TextValue "D5" "580.80"
$ws.Range("E5").Value = "  -4.50%  "
Set-TextValue "D6" "132.79"
$ws.Range("E6").Value = "  -9.04%  "
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.405.25"
$ws.Range("E8").Value = "  -4.24%  "
Set-TextValue "D9" "0.481"
$ws.Range("E9").Value = "  -7.07%  "
$ws.Range("E10").Value = "  -10.18%  "
Set-TextValue "D11" "6.98"
$ws.Range("E11").Value = "  -11.35%  "
Set-TextValue "D12" "0.372"
$ws.Range("E12").Value = "  -10.34%  "
$ws.Range("D13").Value = "3.984.87"
$ws.Range("E13").Value = "  -4.21%  "
Set-TextValue "D14" "0.0000176"
$ws.Range("E14").Value = "  -10.56%  "
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.395.88"
$ws.Range("E16").Value = "  -4.43%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D17" "25.95"
$ws.Range("E17").Value = "  -11.22%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "64.654.02"
$ws.Range("E18").Value = "  -3.17%  "
Set-TextValue "D19" "9.34"
$ws.Range("E19").Value = "  -15.51%  "
Set-TextValue "D20" "5.66"
$ws.Range("E20").Value = "  -9.29%  "
Set-TextValue "D21" "13.38"
$ws.Range("E21").Value = "  -9.48%  "
Set-TextValue "D22" "378.71"
$ws.Range("E22").Value = "  -11.57%  "
$ws.Range("E23").Value = "  +0.04%  "
Set-TextValue "D24" "0.538"
$ws.Range("E24").Value = "  -10.56%  "
Set-TextValue "D25" "71.60"
$ws.Range("E25").Value = "  -7.95%  "
$ws.Range("D26").Value = "3.543.89"
$ws.Range("E26").Value = "  -4.20%  "
Set-TextValue "D27" "0.0000103"
$ws.Range("E27").Value = "  -11.70%  "
Set-TextValue "D28" "1.00"
$ws.Range("E28").Value = "  -0.36%  "
Set-TextValue "D29" "7.14"
$ws.Range("E29").Value = "  -11.93%  "
Set-TextValue "D30" "2.17"
$ws.Range("E30").Value = "  -12.67%  "
Set-TextValue "D31" "7.90"
$ws.Range("E31").Value = "  -12.93%  "
$ws.Range("D32").Value = "3.426.44"
$ws.Range("E32").Value = "  -3.83%  "
$ws.Range("E33").Value = "  -0.02%  "
Set-TextValue "D34" "0.142"
$ws.Range("E34").Value = "  -9.38%  "
Set-TextValue "D35" "22.76"
$ws.Range("E35").Value = "  -7.24%  "
Set-TextValue "D36" "169.92"
$ws.Range("E36").Value = "  -4.10%  "
Set-TextValue "D37" "1.18"
$ws.Range("E37").Value = "  -13.82%  "
Set-TextValue "D38" "6.57"
$ws.Range("E38").Value = "  -15.08%  "
Set-TextValue "D39" "1.43"
$ws.Range("E39").Value = "  -13.27%  "
Set-TextValue "D40" "4.56"
$ws.Range("E40").Value = "  -14.24%  "
Set-TextValue "D41" "0.0754"
$ws.Range("E41").Value = "  -9.45%  "
Set-TextValue "D42" "0.799"
$ws.Range("E42").Value = "  -7.83%  "
$ws.Range("E43").Value = "  +0.12%  "
Set-TextValue "D44" "41.89"
$ws.Range("E44").Value = "  -8.28%  "
Set-TextValue "D45" "4.23"
$ws.Range("E45").Value = "  -16.23%  "
Set-TextValue "D46" "1.59"
$ws.Range("E46").Value = "  -11.85%  "
Set-TextValue "D47" "1.10"
$ws.Range("E47").Value = "  -3.24%  "
Set-TextValue "D48" "22.09"
$ws.Range("E48").Value = "  -7.31%  "
Set-TextValue "D49" "6.45"
$ws.Range("E49").Value = "  -9.99%  "
$ws.Range("D50").Value = "2.194.49"
$ws.Range("E50").Value = "  -6.08%  "
Set-TextValue "D51" "1.94"
$ws.Range("E51").Value = "  -19.89%  "
